$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<Name>_old" -> "<Name>_FV2210" (cols A-J),
#    "diff" stays as-is (col K), "<Name>_new" -> "<Name>_FV2304" (cols L-U)
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value()
    if ($old -like "*_old") {
        $cell.Value = $old.Substring(0, $old.Length - 4) + "_FV2210"
    } elseif ($old -like "*_new") {
        $cell.Value = $old.Substring(0, $old.Length - 4) + "_FV2304"
    }
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into an Excel Table ("Table1") without disturbing
#    the existing header-row formatting / styles.xml (Excel auto-captures a
#    dxf from whatever formatting sits on the header row at the moment
#    ListObjects.Add runs, so stash + restore it around the call).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$stash = $ws.Range("A100:U100")
$headerRange.Copy($stash)
$headerRange.Style = "Normal"

$dataRange = $ws.Range("A1:U66")
$lo = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$stash.Copy()
$headerRange.PasteSpecial(-4122)
$stash.Clear()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
